$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value edits: "paint" the diagonal path on the map ---

# Rows 2 and 3: columns D through AG go from 0 to 100
$ws.Range("D2:AG3").Value = 100

# Diagonal staircase of 110s through columns K/S/AA (rows 8-11 and 17-20)
$ws.Range("K8:K11").Value = 110
$ws.Range("S8:S11").Value = 110
$ws.Range("AA8:AA11").Value = 110

$ws.Range("K17:K20").Value = 110
$ws.Range("S17:S20").Value = 110
$ws.Range("AA17:AA20").Value = 110

# ... and through columns G/O/W (rows 12-15)
$ws.Range("G12:G15").Value = 110
$ws.Range("O12:O15").Value = 110
$ws.Range("W12:W15").Value = 110

# --- Remove the (now redundant) third conditional-formatting rule set ---
# which covered T3:T9, B2:AC2 and AE2:AG2 - it gets cleaned up/merged away,
# leaving just the two remaining rule blocks.
$target = $ws.Range("T3:T9").Address()
$i = 1
while ($i -le $ws.Cells.FormatConditions.Count) {
    $fc = $ws.Cells.FormatConditions.Item($i)
    if ($fc.AppliesTo.Address() -eq $target) {
        $fc.Delete()
    } else {
        $i = $i + 1
    }
}

# --- Update the active selection on the sheet ---
$ws.Range("AF12").Select()
